$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure header for column E (Description) - already present, leave as is.

# For each data row, move the existing Description (column D) value to column E,
# and set the new Booking.com Price value in column D.

$updates = @(
    @{ Row = 2; Price = "92.88 USD" },
    @{ Row = 3; Price = "391.64 USD" },
    @{ Row = 4; Price = "87.82 USD" },
    @{ Row = 5; Price = "88.29 USD" },
    @{ Row = 6; Price = "109 USD" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $desc = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = $desc
    $ws.Cells.Item($r, 4).Value = $u.Price
}

$wb.Save()
